$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing A,B,C
# columns (and their custom widths) to B,C,D, preserving their exact
# widths/content, and leaves a fresh column A for the new "ID COMPANY" data.
$ws.Columns.Item(1).Insert()

# --- Column widths ---
# Column A (new) ~ 14.5703125
$ws.Columns.Item(1).ColumnWidth = 13.6667
# Column E (new) ~ 20.140625
$ws.Columns.Item(5).ColumnWidth = 19.3333

# --- Header row ---
$ws.Range("A1").Value2 = "ID COMPANY"
$ws.Range("D1").Value2 = "CATEGORY ID"
$ws.Range("E1").Value2 = "SUB CATEGORY ID"
$ws.Range("F1").Value2 = "TARGET ID"

# --- Row 2 ---
$ws.Range("A2").Value2 = 2
$ws.Range("D2").Value2 = 1
$ws.Range("E2").Value2 = 5
$ws.Range("F2").Value2 = 1.3

# --- Row 3 ---
$ws.Range("A3").Value2 = 5
$ws.Range("D3").Value2 = 3
$ws.Range("E3").Value2 = 6
$ws.Range("F3").Value2 = 4.8

# --- Selection / view state ---
$ws.Range("E8").Select()
